# Generate Report for Handback
# - The "Ready for handoff" status text becomes "Handback transform failed"
#   everywhere it is used (Overview!E3/F3, zh-cn!C3, de-de!C3).
# - The "Error Detail" column (P) on the zh-cn and de-de sheets gets a
#   diagnostic message for the 6a1d2bf3... row (row 3), explaining the
#   handback/handoff file name mismatch.
# - The Error Detail column is widened to fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Status text update - every cell that showed "Ready for handoff" now shows
# "Handback transform failed" (Overview!E3, Overview!F3, zh-cn!C3, de-de!C3).
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# New Error Detail messages for the 6a1d2bf3... row on each locale sheet.
$zhcn.Range("P3").Value = "Handback file name: nyip1d2p.042 is different with handoff file name: 6a1d2bf3-81d1-49d0-b5a4-179e3295a157.a15ecce4fe6c4da313acff211b68f0d1c82a452b.zh-cn."
$dede.Range("P3").Value = "Handback file name: nyip1d2p.042 is different with handoff file name: 6a1d2bf3-81d1-49d0-b5a4-179e3295a157.a15ecce4fe6c4da313acff211b68f0d1c82a452b.de-de."

# Widen the Error Detail column (P = column 16) on both locale sheets.
# (ColumnWidth 39.17 "characters" renders as the OOXML width="40" units
# used by the stored column definition.)
$zhcn.Columns.Item(16).ColumnWidth = 39.17
$dede.Columns.Item(16).ColumnWidth = 39.17
